$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col4a6"
$ws.Range("C2").Value = "Cd93"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03106633333333333
$ws.Range("H2").Value = 0.093199
$ws.Range("I2").Value = 0.009330498513261214
$ws.Range("J2").Value = 0.009330498513261214
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 124.722578
$ws.Range("N2").Value = 374.167734
$ws.Range("O2").Value = 0.9767311432246923
$ws.Range("P2").Value = 0.9767311432246923
$ws.Range("Q2").Value = 3.874673182340667
$ws.Range("R2").Value = 34.872058641066
$ws.Range("S2").Value = 0.009113388479713918
$ws.Range("T2").Value = 0.009113388479713918

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col4a6"
$ws.Range("C3").Value = "Cd93"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03106633333333333
$ws.Range("H3").Value = 0.093199
$ws.Range("I3").Value = 0.009330498513261214
$ws.Range("J3").Value = 0.009330498513261214
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.134712
$ws.Range("N3").Value = 0.404136
$ws.Range("O3").Value = 0.001054960600366076
$ws.Range("P3").Value = 0.001054960600366076
$ws.Range("Q3").Value = 0.004185007896
$ws.Range("R3").Value = 0.037665071064
$ws.Range("S3").Value = 0.0000098433083132648301815
$ws.Range("T3").Value = 0.0000098433083132648301815

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col4a6"
$ws.Range("C4").Value = "Cd93"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03106633333333333
$ws.Range("H4").Value = 0.093199
$ws.Range("I4").Value = 0.009330498513261214
$ws.Range("J4").Value = 0.009330498513261214
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509735
$ws.Range("O4").Value = 0.02221389617494163
$ws.Range("P4").Value = 0.02221389617494163
$ws.Range("Q4").Value = 0.08812208802944445
$ws.Range("R4").Value = 0.7930987922649999
$ws.Range("S4").Value = 0.0002072667252340318
$ws.Range("T4").Value = 0.0002072667252340318

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col4a6"
$ws.Range("C5").Value = "Cd93"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9795656666666667
$ws.Range("H5").Value = 2.938697
$ws.Range("I5").Value = 0.2942038861943281
$ws.Range("J5").Value = 0.2942038861943281
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 124.722578
$ws.Range("N5").Value = 374.167734
$ws.Range("O5").Value = 0.9767311432246923
$ws.Range("P5").Value = 0.9767311432246923
$ws.Range("Q5").Value = 122.1739552669553
$ws.Range("R5").Value = 1099.565597402598
$ws.Range("S5").Value = 0.2873580981037334
$ws.Range("T5").Value = 0.2873580981037334

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col4a6"
$ws.Range("C6").Value = "Cd93"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9795656666666667
$ws.Range("H6").Value = 2.938697
$ws.Range("I6").Value = 0.2942038861943281
$ws.Range("J6").Value = 0.2942038861943281
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.134712
$ws.Range("N6").Value = 0.404136
$ws.Range("O6").Value = 0.001054960600366076
$ws.Range("P6").Value = 0.001054960600366076
$ws.Range("Q6").Value = 0.131959250088
$ws.Range("R6").Value = 1.187633250792
$ws.Range("S6").Value = 0.0003103735084096011
$ws.Range("T6").Value = 0.0003103735084096011

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col4a6"
$ws.Range("C7").Value = "Cd93"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9795656666666667
$ws.Range("H7").Value = 2.938697
$ws.Range("I7").Value = 0.2942038861943281
$ws.Range("J7").Value = 0.2942038861943281
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509735
$ws.Range("O7").Value = 0.02221389617494163
$ws.Range("P7").Value = 0.02221389617494163
$ws.Range("Q7").Value = 2.778614746143889
$ws.Range("R7").Value = 25.007532715295
$ws.Range("S7").Value = 0.006535414582185148
$ws.Range("T7").Value = 0.006535414582185148

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col4a6"
$ws.Range("C8").Value = "Cd93"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.318915
$ws.Range("H8").Value = 6.956745
$ws.Range("I8").Value = 0.6964656152924107
$ws.Range("J8").Value = 0.6964656152924107
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 124.722578
$ws.Range("N8").Value = 374.167734
$ws.Range("O8").Value = 0.9767311432246923
$ws.Range("P8").Value = 0.9767311432246923
$ws.Range("Q8").Value = 289.22105696287
$ws.Range("R8").Value = 2602.98951266583
$ws.Range("S8").Value = 0.6802596566412451
$ws.Range("T8").Value = 0.6802596566412451

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col4a6"
$ws.Range("C9").Value = "Cd93"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.318915
$ws.Range("H9").Value = 6.956745
$ws.Range("I9").Value = 0.6964656152924107
$ws.Range("J9").Value = 0.6964656152924107
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.134712
$ws.Range("N9").Value = 0.404136
$ws.Range("O9").Value = 0.001054960600366076
$ws.Range("P9").Value = 0.001054960600366076
$ws.Range("Q9").Value = 0.31238567748
$ws.Range("R9").Value = 2.81147109732
$ws.Range("S9").Value = 0.0007347437836432102
$ws.Range("T9").Value = 0.0007347437836432102

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col4a6"
$ws.Range("C10").Value = "Cd93"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.318915
$ws.Range("H10").Value = 6.956745
$ws.Range("I10").Value = 0.6964656152924107
$ws.Range("J10").Value = 0.6964656152924107
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509735
$ws.Range("O10").Value = 0.02221389617494163
$ws.Range("P10").Value = 0.02221389617494163
$ws.Range("Q10").Value = 6.577784045841667
$ws.Range("R10").Value = 59.20005641257499
$ws.Range("S10").Value = 0.01547121486752245
$ws.Range("T10").Value = 0.01547121486752245

